# Generate Report for Handoff
# Refresh the localization-status report: the handoff package moved from
# GUID 5b18ebab-4989-4f77-9d74-9f18a4b8679b to fbfa8d07-4749-4eb6-bf91-fa62043e61c5
# (new handoff run), the xlf build hash moved from 147a98fc5c73076d2c19aea123d8f5a90a4bb092
# to 6c76c397a99cc3bac9c1d7958a1dedcd08537545, and the handoff timestamps advanced.

$wb = $excel.ActiveWorkbook

$oldGuid = "5b18ebab-4989-4f77-9d74-9f18a4b8679b"
$newGuid = "fbfa8d07-4749-4eb6-bf91-fa62043e61c5"
$oldHash = "147a98fc5c73076d2c19aea123d8f5a90a4bb092"
$newHash = "6c76c397a99cc3bac9c1d7958a1dedcd08537545"

$newMdName    = $newGuid + ".md"
$newZhCnName  = $newGuid + "." + $newHash + ".zh-cn.xlf"
$newDeDeName  = $newGuid + "." + $newHash + ".de-de.xlf"

$newHandoffDate       = "2016-03-21 04:59:18"
$newZhCnHandoffDate   = "2016-03-21 04:59:09"
$newDeDeHandoffDate   = "2016-03-21 04:59:18"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newHandoffDate

foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = $newMdName
}

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("D2").Value = $newZhCnName
$wsZhCn.Range("E2").Value = $newZhCnHandoffDate

foreach ($h in $wsZhCn.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newMdName
    } elseif ($addr -eq '$D$2') {
        $h.TextToDisplay = $newZhCnName
    }
}

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("D2").Value = $newDeDeName
$wsDeDe.Range("E2").Value = $newDeDeHandoffDate

foreach ($h in $wsDeDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newMdName
    } elseif ($addr -eq '$D$2') {
        $h.TextToDisplay = $newDeDeName
    }
}
